# Scheduled-runner price/profit refresh across the Leve worksheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) for the rows whose market data changed.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 39
$ws.Range("H39").Value = 456.53845
$ws.Range("I39").Value = 228
$ws.Range("J39").Value = 652.4286
$ws.Range("K39").Value = 684
$ws.Range("L39").Value = 1957.2858
$ws.Range("M39").Value = -388
$ws.Range("N39").Value = -2549.2858

# row 138
$ws.Range("H138").Value = 2110.65
$ws.Range("I138").Value = 883.25
$ws.Range("J138").Value = 2278.0227
$ws.Range("K138").Value = 2649.75
$ws.Range("L138").Value = 6834.0681
$ws.Range("M138").Value = 2490.25
$ws.Range("N138").Value = -17114.0681

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 5
$ws.Range("H5").Value = 138.83333
$ws.Range("I5").Value = 128.25
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 128.25
$ws.Range("L5").Value = 160
$ws.Range("M5").Value = -16.25
$ws.Range("N5").Value = -384

# row 61
$ws.Range("H61").Value = 3069.7334
$ws.Range("I61").Value = 2620.0588
$ws.Range("J61").Value = 3657.7693
$ws.Range("K61").Value = 2620.0588
$ws.Range("L61").Value = 3657.7693
$ws.Range("M61").Value = -2408.0588
$ws.Range("N61").Value = -4081.7693

# row 63
$ws.Range("H63").Value = 5436.077
$ws.Range("I63").Value = 3241.8572
$ws.Range("J63").Value = 7996
$ws.Range("K63").Value = 3241.8572
$ws.Range("L63").Value = 7996
$ws.Range("M63").Value = -2555.8572
$ws.Range("N63").Value = -9368

# row 66
$ws.Range("H66").Value = 5436.077
$ws.Range("I66").Value = 3241.8572
$ws.Range("J66").Value = 7996
$ws.Range("K66").Value = 16209.286
$ws.Range("L66").Value = 39980
$ws.Range("M66").Value = -12777.286
$ws.Range("N66").Value = -46844

# row 74
$ws.Range("H74").Value = 3236.4285
$ws.Range("I74").Value = 2980.2727
$ws.Range("J74").Value = 4175.6665
$ws.Range("K74").Value = 2980.2727
$ws.Range("L74").Value = 4175.6665
$ws.Range("M74").Value = -2106.2727
$ws.Range("N74").Value = -5923.6665

# row 77
$ws.Range("H77").Value = 3236.4285
$ws.Range("I77").Value = 2980.2727
$ws.Range("J77").Value = 4175.6665
$ws.Range("K77").Value = 14901.3635
$ws.Range("L77").Value = 20878.3325
$ws.Range("M77").Value = -10533.3635
$ws.Range("N77").Value = -29614.3325

# row 136
$ws.Range("H136").Value = 3069.7334
$ws.Range("I136").Value = 2620.0588
$ws.Range("J136").Value = 3657.7693
$ws.Range("K136").Value = 7860.176399999999
$ws.Range("L136").Value = 10973.3079
$ws.Range("M136").Value = -5310.176399999999
$ws.Range("N136").Value = -16073.3079

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 138.83333
$ws.Range("I4").Value = 128.25
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 128.25
$ws.Range("L4").Value = 160
$ws.Range("M4").Value = -13.25
$ws.Range("N4").Value = -390

# row 15
$ws.Range("H15").Value = 50000000
$ws.Range("I15").Value = 50000000
$ws.Range("K15").Value = 50000000
$ws.Range("L15").Value = -49999773

# row 35
$ws.Range("H35").Value = 16964.166
$ws.Range("J35").Value = 19357
$ws.Range("L35").Value = 19357
$ws.Range("N35").Value = -19977

# row 82
$ws.Range("H82").Value = 13907.066
$ws.Range("J82").Value = 21249
$ws.Range("L82").Value = 21249
$ws.Range("N82").Value = -22015

# row 85
$ws.Range("H85").Value = 13907.066
$ws.Range("J85").Value = 21249
$ws.Range("L85").Value = 21249
$ws.Range("M85").Value = -4190.2856
$ws.Range("N85").Value = -23901

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 333.25
$ws.Range("I22").Value = 269.16666
$ws.Range("J22").Value = 397.33334
$ws.Range("K22").Value = 269.16666
$ws.Range("L22").Value = 397.33334
$ws.Range("M22").Value = 80.83334000000002
$ws.Range("N22").Value = -1097.33334

# row 94
$ws.Range("H94").Value = 1196.9231
$ws.Range("I94").Value = 826.6667
$ws.Range("J94").Value = 1308
$ws.Range("K94").Value = 826.6667
$ws.Range("L94").Value = 1308
$ws.Range("M94").Value = -375.6667
$ws.Range("N94").Value = -2210

# row 134
$ws.Range("H134").Value = 1573.6428
$ws.Range("I134").Value = 1367.3636
$ws.Range("J134").Value = 2330
$ws.Range("K134").Value = 4102.0908
$ws.Range("L134").Value = 6990
$ws.Range("M134").Value = -1567.0908
$ws.Range("N134").Value = -12060

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 98
$ws.Range("H98").Value = 266.25
$ws.Range("I98").Value = 230
$ws.Range("J98").Value = 326.66666
$ws.Range("K98").Value = 690
$ws.Range("L98").Value = 979.9999799999999
$ws.Range("M98").Value = 808
$ws.Range("N98").Value = -3975.99998

# row 122
$ws.Range("H122").Value = 13011
$ws.Range("I122").Value = 348.33334
$ws.Range("J122").Value = 50999
$ws.Range("K122").Value = 3135.00006
$ws.Range("L122").Value = 458991
$ws.Range("M122").Value = -685.0000600000003
$ws.Range("N122").Value = -463891

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 79.63636
$ws.Range("I2").Value = 43.4
$ws.Range("J2").Value = 109.833336
$ws.Range("K2").Value = 43.4
$ws.Range("L2").Value = 109.833336
$ws.Range("M2").Value = 69.59999999999999
$ws.Range("N2").Value = -335.833336

# row 18
$ws.Range("H18").Value = 8500
$ws.Range("J18").Value = 8500
$ws.Range("L18").Value = 8500
$ws.Range("N18").Value = -9086

# row 39
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064

# row 43
$ws.Range("H43").Value = 14179.4
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 14179.4
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 14179.4
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -14481.4

# row 80
$ws.Range("H80").Value = 33962824
$ws.Range("I80").Value = 42419700
$ws.Range("J80").Value = 135333.33
$ws.Range("K80").Value = 42419700
$ws.Range("L80").Value = 135333.33
$ws.Range("M80").Value = -42418702
$ws.Range("N80").Value = -137329.33

# row 83
$ws.Range("H83").Value = 33962824
$ws.Range("I83").Value = 42419700
$ws.Range("J83").Value = 135333.33
$ws.Range("K83").Value = 212098500
$ws.Range("L83").Value = 676666.6499999999
$ws.Range("M83").Value = -212093508
$ws.Range("N83").Value = -686650.6499999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 11729.667
$ws.Range("I22").Value = 1233.3334
$ws.Range("J22").Value = 13828.934
$ws.Range("K22").Value = 1233.3334
$ws.Range("L22").Value = 13828.934
$ws.Range("M22").Value = -938.3334
$ws.Range("N22").Value = -14418.934

# row 27
$ws.Range("H27").Value = 11729.667
$ws.Range("I27").Value = 1233.3334
$ws.Range("J27").Value = 13828.934
$ws.Range("K27").Value = 1233.3334
$ws.Range("L27").Value = 13828.934
$ws.Range("M27").Value = -1126.3334
$ws.Range("N27").Value = -14042.934

# row 46
$ws.Range("H46").Value = 5849.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5849.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5849.75
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6225.75

# row 132
$ws.Range("H132").Value = 2357.9167
$ws.Range("I132").Value = 1399.4546
$ws.Range("K132").Value = 4198.3638
$ws.Range("M132").Value = -1668.3638

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 3789955.2
$ws.Range("I132").Value = 1945.7037
$ws.Range("K132").Value = 5837.1111
$ws.Range("M132").Value = -3307.1111
